$wb = $excel.ActiveWorkbook

# --- Sheet "summary" ---
$ws1 = $wb.Worksheets.Item("summary")
$ws1.Range("B2").Value = 41.88
$ws1.Range("D2").Value = 298
$ws1.Range("B3").Value = 24.46
$ws1.Range("C3").Value = 14.58
$ws1.Range("D3").Value = 301

# --- Sheet "summary_repository" ---
$ws2 = $wb.Worksheets.Item("summary_repository")
$ws2.Range("C3").Value = 68.88
$ws2.Range("E3").Value = 80
$ws2.Range("C4").Value = 14.18
$ws2.Range("D4").Value = 14.58
$ws2.Range("E4").Value = 219
$ws2.Range("C5").Value = 51.9
$ws2.Range("E5").Value = 82
